$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.827.55"
$ws.Range("E2").Value = "  -5.74%  "

# Row 3
$ws.Range("D3").Value = "2.226.73"
$ws.Range("E3").Value = "  -6.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'312.99"
$ws.Range("E5").Value = "  -1.58%  "

# Row 6
$ws.Range("D6").Value = "'100.86"
$ws.Range("E6").Value = "  -10.98%  "

# Row 7
$ws.Range("D7").Value = "'0.570"
$ws.Range("E7").Value = "  -10.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.564"
$ws.Range("E9").Value = "  -9.81%  "

# Row 10
$ws.Range("D10").Value = "'37.58"
$ws.Range("E10").Value = "  -10.68%  "

# Row 11
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  -8.82%  "

# Row 12
$ws.Range("D12").Value = "'7.68"
$ws.Range("E12").Value = "  -11.32%  "

# Row 13
$ws.Range("E13").Value = "  -4.25%  "

# Row 14
$ws.Range("D14").Value = "'0.889"
$ws.Range("E14").Value = "  -11.35%  "

# Row 15
$ws.Range("D15").Value = "2.558.64"
$ws.Range("E15").Value = "  -7.29%  "

# Row 16
$ws.Range("D16").Value = "'13.89"
$ws.Range("E16").Value = "  -12.26%  "

# Row 17
$ws.Range("D17").Value = "2.195.06"
$ws.Range("E17").Value = "  -8.60%  "

# Row 18
$ws.Range("D18").Value = "42.553.53"
$ws.Range("E18").Value = "  -6.33%  "

# Row 19
$ws.Range("D19").Value = "'14.50"
$ws.Range("E19").Value = "  +8.28%  "

# Row 20
$ws.Range("D20").Value = "'6.68"
$ws.Range("E20").Value = "  -10.81%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -11.28%  "

# Row 22
$ws.Range("D22").Value = "'3.28"
$ws.Range("E22").Value = "  -7.60%  "

# Row 23
$ws.Range("D23").Value = "'65.35"
$ws.Range("E23").Value = "  -12.43%  "

# Row 24
$ws.Range("D24").Value = "'236.58"
$ws.Range("E24").Value = "  -10.45%  "

# Row 25
$ws.Range("D25").Value = "'2.16"
$ws.Range("E25").Value = "  -8.19%  "

# Row 26
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  +0.47%  "

# Row 27
$ws.Range("D27").Value = "'10.35"
$ws.Range("E27").Value = "  -8.44%  "

# Row 28
$ws.Range("D28").Value = "'6.67"
$ws.Range("E28").Value = "  -13.54%  "

# Row 29
$ws.Range("D29").Value = "'2.16"
$ws.Range("E29").Value = "  -8.43%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0905"
$ws.Range("E30").Value = "  -6.54%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.77"
$ws.Range("E31").Value = "  -8.51%  "

# Row 32
$ws.Range("D32").Value = "'34.69"
$ws.Range("E32").Value = "  -11.80%  "

# Row 33
$ws.Range("D33").Value = "'159.74"
$ws.Range("E33").Value = "  -7.47%  "

# Row 34
$ws.Range("D34").Value = "'2.74"
$ws.Range("E34").Value = "  -6.93%  "

# Row 35
$ws.Range("D35").Value = "'3.17"
$ws.Range("E35").Value = "  +2.90%  "

# Row 36
$ws.Range("D36").Value = "'0.122"
$ws.Range("E36").Value = "  -7.92%  "

# Row 37
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "  +7.44%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.37"
$ws.Range("E38").Value = "  -11.25%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.104"
$ws.Range("E39").Value = "  -12.08%  "

# Row 40
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  -10.96%  "

# Row 41
$ws.Range("D41").Value = "'0.0323"
$ws.Range("E41").Value = "  -10.95%  "

# Row 42
$ws.Range("E42").Value = "  -0.11%  "

# Row 43
$ws.Range("D43").Value = "1.811.39"
$ws.Range("E43").Value = "  +8.35%  "

# Row 44
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "'91.59"
$ws.Range("E44").Value = "  -11.06%  "

# Row 45
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'12.34"
$ws.Range("E45").Value = "  -9.77%  "

# Row 46
$ws.Range("D46").Value = "'0.209"
$ws.Range("E46").Value = "  -12.73%  "

# Row 47
$ws.Range("D47").Value = "'62.00"
$ws.Range("E47").Value = "  -13.18%  "

# Row 48
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'77.95"
$ws.Range("E48").Value = "  -11.01%  "

# Row 49
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'5.37"
$ws.Range("E49").Value = "  -6.83%  "

# Row 50
$ws.Range("D50").Value = "'8.56"
$ws.Range("E50").Value = "  -9.66%  "

# Row 51
$ws.Range("D51").Value = "'102.88"
$ws.Range("E51").Value = "  -11.26%  "
